$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark -------------------------------
# In the original document the _GoBack bookmark sits at the very end of
# the "Regroup on Wednesday" paragraph. The edit moves it so that it
# instead sits right after the "Product_suppliers" run (immediately
# before the spellEnd proofing mark) in the
# "Add/Edit product, suppliers, product_suppliers" bullet's sub-item.
# Re-adding a bookmark with the same name relocates the existing one
# (bookmark names are unique), so this single Add() both removes it
# from its old spot and places it in the new one.
$p = $d.Paragraphs(8)
$paraText = $p.Range.Text
$target = "Product_suppliers"
$idx = $paraText.IndexOf($target)
$pos = $p.Range.Start + $idx + $target.Length
$bookmarkRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 2. Drop the "does not need to be displayed" annotations ----------
# Two list items ("Product_suppliers" and "Packages_products_suppliers")
# lose their trailing " //does not need to be displayed" comment text.
$d.Content.Find.Execute(" //does not need to be displayed", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
